# ---------------------------------------------------------------------------
# This script reproduces the commit:
#   "se deja listo fechas y put" -- adds an " {apellido}" run right after the
#   "{nombre}" run (with identical run formatting), and adds a "Texto" run
#   right after the "fechaInicio" run (with identical run formatting).
#
# Both insertions need to land in the saved document as their OWN <w:r>
# element (not merged into the neighboring run's <w:t>), even though their
# run formatting (<w:rPr>) is byte-for-byte identical to the run that
# precedes them. Word's "type text in place" operations (Range.Text = ...,
# Range.InsertAfter(...), Find.Execute(..., Replace:=...)) immediately fold
# newly typed text into the neighboring run whenever the effective
# formatting matches, so a literal in-place insert would merge the two
# pieces of text into a single run. Copying a Range's .FormattedText onto
# another Range, however, preserves run identity, so we use a tiny
# off-document scratch paragraph to "clone" the desired formatting, edit the
# clone's text, and then paste that clone's FormattedText right after the
# target run; the scratch paragraph is removed afterward.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Insert-RunAfter($findText, $newText) {
    # Locate the run whose text we want to extend with a sibling run.
    $target = $d.Content
    $target.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Create a scratch paragraph at the very end of the document body; we
    # will use it to build a clone of the target run's formatting without
    # disturbing the target itself.
    $end = $d.Content
    $end.Collapse(0)
    $end.InsertParagraphAfter()
    $scratchIndex = $d.Paragraphs.Count
    $scratchPara = $d.Paragraphs($scratchIndex)
    $scratch = $scratchPara.Range

    # Clone the found run's formatted text (this keeps it as its own run
    # rather than merging it with whatever else is in the scratch paragraph).
    $scratch.FormattedText = $target.FormattedText

    # Restrict to just the cloned run (exclude the paragraph mark) and
    # replace its text with the text we actually want to insert.
    $clone = $scratchPara.Range
    $clone.MoveEnd(1, -1)
    $clone.Text = $newText

    # Re-locate the target run (indexes are stable; Range objects captured
    # earlier still point at the right place) and paste the clone's
    # formatted text right after it -- this adds a new sibling run with the
    # same run formatting instead of merging into the target run.
    $target.Collapse(0)
    $target.FormattedText = $clone.FormattedText

    # Remove the temporary scratch paragraph we used to build the clone.
    $lastIndex = $d.Paragraphs.Count
    $d.Paragraphs($lastIndex).Range.Delete()
}

Insert-RunAfter "{nombre}" " {apellido}"
Insert-RunAfter "fechaInicio" "Texto"
